$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete column F ("Department") first so that column B's index is unaffected.
$ws.Range("F1").EntireColumn.Delete()
# Delete column B ("Program Title").
$ws.Range("B1").EntireColumn.Delete()

# Column G ("Provided Funds (INR)") ends up selected after the edit.
$ws.Columns("G:G").Select()
